$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativação: date text updated (row 8, B/C). Written via a text-literal
# formula + paste-as-values so Excel keeps it as a plain text string rather
# than auto-converting the dd/mm/yyyy-looking text into a date serial. ---
foreach ($addr in @("B8","C8")) {
  $cell = $ws.Range($addr)
  $cell.Formula = "=""01/01/2022"""
  $cell.Copy()
  $cell.PasteSpecial(-4163)   # xlPasteValues
}
$excel.CutCopyMode = $false

# --- Docentes responsáveis: first teacher changed (row 13, not shifted) ---
$ws.Range("B13").Value2 = "5840897 - Clodoaldo Saron"
$ws.Range("C13").Value2 = "5840897 - Clodoaldo Saron"

# --- Insert a new row after row 13 to hold the second teacher's name ---
# (shifts old rows 14-24 down to 15-25)
$ws.Rows.Item(14).EntireRow.Insert()
$ws.Range("B14").Value2 = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C14").Value2 = "1033242 - Fábio Herbst Florenzano"

# --- Método: evaluation method text changed (old row 19, now row 20) ---
$ws.Range("B20").Value2 = "A avaliação será feita por meio de Provas Escritas, Estudos de Casos e Desenvolvimento de Projetos, sendo necessário utilizar pelo menos dois critérios de avaliação diferentes."
$ws.Range("C20").Value2 = "A avaliação será feita por meio de Provas Escritas, Estudos de Casos e Desenvolvimento de Projetos, sendo necessário utilizar pelo menos dois critérios de avaliação diferentes."

# --- Critério: grade formula text changed (old row 20, now row 21) ---
$ws.Range("B21").Value2 = "A Nota final (NF) será calculada da seguinte maneira: NF = (P+EC+Projetos)/3"
$ws.Range("C21").Value2 = "A Nota final (NF) será calculada da seguinte maneira: NF = (P+EC+Projetos)/3"

# --- Norma de recuperação: recovery norm text changed (old row 21, now row 22) ---
$ws.Range("B22").Value2 = "Não consta recuperação"
$ws.Range("C22").Value2 = "Não consta recuperação"

# --- Bibliografia: bibliography text changed (old row 22, now row 23) ---
$ws.Range("B23").Value2 = "1. J. Margolis. Engineering Plastics Handbook. McGraw-Hill Professional, 2005. 2. Nigel Mills. Plastics - Microstructure and Engineering Applications. Butterworth-Heineman, 2005. 3. Walter Michaeli, TEcnologia dos Plasticos. Ed. Blucher 4. Hélio Wiebeck, Júlio Harada. Plásticos de Engenharia - Tecnologia e Aplicações. São Paulo: Editora Artliber, 2005. 5. E. B. Mano, L. C. Mendes. Identificação de Plásticos, Borrachas e Fibras. São Paulo: Editora Edgard Blucher, 2000. 6. Marcelo Rabello. Aditivação de Polímeros. São Paulo: Editora Artliber, 2004. 7. Jan C.J. Bart. Additives in Polymers. New York: John Wiley & Sons, 2005. 8. Marino Xanthos. Functional Fillers for Plastics. Wiley-VCH Verlag GmbH, 2005. 9. Silvio Manrich. Processamento de Termoplásticos. Editora Artliber, 2005. 10. G.H. Michler, F.J. Baltá-Calleja. Mechanical Properties of Polymers Based on Nanostructure and Morphology. Boca Raton: CRC Press, 2005. 11. A. M. Piva, H. Wiebeck. Reciclagem do P. São Paulo: Editora Artliber"". Manas Chanda, ,Salil K. Roy  Plastics Fabrication and Recycling"
$ws.Range("C23").Value2 = "1. J. Margolis. Engineering Plastics Handbook. McGraw-Hill Professional, 2005. 2. Nigel Mills. Plastics - Microstructure and Engineering Applications. Butterworth-Heineman, 2005. 3. Walter Michaeli, TEcnologia dos Plasticos. Ed. Blucher 4. Hélio Wiebeck, Júlio Harada. Plásticos de Engenharia - Tecnologia e Aplicações. São Paulo: Editora Artliber, 2005. 5. E. B. Mano, L. C. Mendes. Identificação de Plásticos, Borrachas e Fibras. São Paulo: Editora Edgard Blucher, 2000. 6. Marcelo Rabello. Aditivação de Polímeros. São Paulo: Editora Artliber, 2004. 7. Jan C.J. Bart. Additives in Polymers. New York: John Wiley & Sons, 2005. 8. Marino Xanthos. Functional Fillers for Plastics. Wiley-VCH Verlag GmbH, 2005. 9. Silvio Manrich. Processamento de Termoplásticos. Editora Artliber, 2005. 10. G.H. Michler, F.J. Baltá-Calleja. Mechanical Properties of Polymers Based on Nanostructure and Morphology. Boca Raton: CRC Press, 2005. 11. A. M. Piva, H. Wiebeck. Reciclagem do P. São Paulo: Editora Artliber"". Manas Chanda, ,Salil K. Roy  Plastics Fabrication and Recycling"
